$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Variable" column), shifting
# Type/Index/Original/Translation one column to the right.
$ws.Range("B1").EntireColumn.Insert()

$ws.Range("B1").Value = "Variable"
$ws.Range("B2").Value = "c1"
$ws.Range("B3").Value = "c1"
$ws.Range("B4").Value = "c1"
$ws.Range("B5").Value = "c1"

$ws.Range("B6").Select() | Out-Null
